$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold number-looking text (e.g. "312.54", "27.108.49")
# that must stay stored as text, not get auto-converted to numbers by Excel's
# smart entry. Mark the range as Text before writing, then restore the
# original (unstyled) look so no stray number-format / quote-prefix style
# gets attached to the cells.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.108.49"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.831.22"
$ws.Range("E3").Value = "  +0.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "312.54"
$ws.Range("E5").Value = "  -0.17%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.10%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4626"
$ws.Range("E7").Value = "  -1.57%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3703"
$ws.Range("E8").Value = "  +1.06%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07354"
$ws.Range("E9").Value = "  -0.64%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.8732"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.08000"
$ws.Range("E11").Value = "  +4.36%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  -2.18%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.808.46"
$ws.Range("E13").Value = "  -4.74%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.346"
$ws.Range("E14").Value = "  -0.63%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "6.564"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "91.97"
$ws.Range("E16").Value = "  -1.54%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000008885"
$ws.Range("E18").Value = "  +1.71%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.16%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "14.68"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.012.79"
$ws.Range("E21").Value = "  -2.11%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.140"
$ws.Range("E22").Value = "  -2.01%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.022.35"
$ws.Range("E24").Value = "  -3.10%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.54"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "1.833"
$ws.Range("E26").Value = "  -2.78%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  +0.77%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "2.085"
$ws.Range("E28").Value = "  -2.12%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "5.094"
$ws.Range("E29").Value = "  -1.56%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "115.40"
$ws.Range("E30").Value = "  -1.05%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.08872"
$ws.Range("E31").Value = "  -0.71%  "

# Row 32 - HuobiToken
$ws.Range("D32").Value = "2.976"
$ws.Range("E32").Value = "  +1.14%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "0.7343"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "4.447"
$ws.Range("E34").Value = "  -1.57%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -2.29%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "2.459"
$ws.Range("E36").Value = "  -4.55%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "1.076"
$ws.Range("E37").Value = "  -1.42%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  +0.50%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "0.05244"
$ws.Range("E39").Value = "  -1.12%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.948"
$ws.Range("E40").Value = "  +0.43%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "7.156"
$ws.Range("E41").Value = "  -2.54%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "0.5195"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "0.1634"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44 - Frax
$ws.Range("D44").Value = "0.8602"
$ws.Range("E44").Value = "  -14.77%  "

# Row 45 - Aptos
$ws.Range("D45").Value = "8.236"
$ws.Range("E45").Value = "  -1.78%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "0.4845"
$ws.Range("E46").Value = "  -1.22%  "

# Row 47 - was EnergySwap, now PaxDollar
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48 - was PaxDollar, now EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.22"
$ws.Range("E48").Value = "  -1.62%  "

# Row 49 - Quant
$ws.Range("D49").Value = "102.47"
$ws.Range("E49").Value = "  -1.91%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "1.631"
$ws.Range("E50").Value = "  -1.46%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.06232"
$ws.Range("E51").Value = "  -0.74%  "

# Restore the default (unstyled) look on the Price/Volume range so the
# saved file has no stray number-format style attached to these cells,
# matching the original formatting.
$ws.Range("D2:E51").Style = "Normal"
